$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph, then work on the
# bullet paragraphs that immediately follow it (skipping the "Impact" sub-heading).
$achIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "KEY ACHIEVEMENTS AND IMPACT") {
        $achIdx = $i
        break
    }
}

if ($achIdx -eq 0) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# Bullets start two paragraphs after the heading (heading, then "Impact" sub-heading, then bullets).
$bulletStart = $achIdx + 2

# Replace the text of the first four bullet paragraphs with the new accomplishment statements.
# NOTE: do not append a trailing `r — assigning Range.Text with an embedded
# paragraph mark inserts an *extra* empty paragraph instead of just replacing
# the existing paragraph's text.
$d.Paragraphs($bulletStart + 0).Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
$d.Paragraphs($bulletStart + 1).Range.Text = "• Real-time collaboration at national scale"
$d.Paragraphs($bulletStart + 2).Range.Text = "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
$d.Paragraphs($bulletStart + 3).Range.Text = "• 23% conversion rate improvement"

# Delete the now-obsolete 5th and 6th bullet paragraphs (their text has been
# superseded, so the paragraphs themselves are removed).
$d.Paragraphs($bulletStart + 4).Range.Delete()
$d.Paragraphs($bulletStart + 4).Range.Delete()
